$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 847.2963
$ws.Range("J17").Value = 690.38464
$ws.Range("L17").Value = 2071.15392
$ws.Range("N17").Value = -2407.15392

$ws.Range("H21").Value = 18333.334
$ws.Range("J21").Value = 20000
$ws.Range("L21").Value = 20000
$ws.Range("N21").Value = -20936

$ws.Range("H23").Value = 18333.334
$ws.Range("J23").Value = 20000
$ws.Range("L23").Value = 20000
$ws.Range("N23").Value = -20468

$ws.Range("H113").Value = 2699.5
$ws.Range("J113").Value = 1948
$ws.Range("L113").Value = 1948
$ws.Range("N113").Value = -8456

$ws.Range("H132").Value = 2699.3044
$ws.Range("J132").Value = 14999.5
$ws.Range("L132").Value = 44998.5
$ws.Range("N132").Value = -50058.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 40.5
$ws.Range("I5").Value = 33.42857
$ws.Range("K5").Value = 33.42857
$ws.Range("M5").Value = 78.57142999999999

$ws.Range("H58").Value = 3539.5
$ws.Range("I58").Value = 2079
$ws.Range("J58").Value = 5000
$ws.Range("K58").Value = 2079
$ws.Range("L58").Value = 5000
$ws.Range("M58").Value = -1649
$ws.Range("N58").Value = -5860

$ws.Range("H62").Value = 0
$ws.Range("I62").Value = 0
$ws.Range("K62").Value = 0
$ws.Range("M62").ClearContents()

$ws.Range("H65").Value = 0
$ws.Range("I65").Value = 0
$ws.Range("K65").Value = 0
$ws.Range("M65").ClearContents()

$ws.Range("H94").Value = 30000
$ws.Range("J94").Value = 30000
$ws.Range("L94").Value = 30000
$ws.Range("N94").Value = -31802

$ws.Range("H96").Value = 23844
$ws.Range("J96").Value = 23844
$ws.Range("L96").Value = 23844
$ws.Range("N96").Value = -29336

$ws.Range("H102").Value = 11113240
$ws.Range("I102").Value = 18519866
$ws.Range("K102").Value = 18519866
$ws.Range("M102").Value = -18518244

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 40.5
$ws.Range("I4").Value = 33.42857
$ws.Range("K4").Value = 33.42857
$ws.Range("M4").Value = 81.57142999999999

$ws.Range("H20").Value = 4142.0625
$ws.Range("I20").Value = 3296.1
$ws.Range("J20").Value = 5552
$ws.Range("K20").Value = 3296.1
$ws.Range("L20").Value = 5552
$ws.Range("M20").Value = -3049.1
$ws.Range("N20").Value = -6046

$ws.Range("H31").Value = 6000
$ws.Range("I31").Value = 6000
$ws.Range("J31").Value = 0
$ws.Range("K31").Value = 6000
$ws.Range("L31").ClearContents()
$ws.Range("M31").Value = -5748
$ws.Range("N31").Value = 0

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H54").Value = 13000
$ws.Range("J54").Value = 13000
$ws.Range("L54").Value = 13000
$ws.Range("N54").Value = -14316

$ws.Range("H55").Value = 0
$ws.Range("I55").Value = 0
$ws.Range("K55").Value = 0
$ws.Range("M55").ClearContents()

$ws.Range("H103").Value = 0
$ws.Range("I103").Value = 0
$ws.Range("K103").Value = 0
$ws.Range("M103").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H23").Value = 333345.34
$ws.Range("I23").Value = 18.5
$ws.Range("K23").Value = 55.5
$ws.Range("M23").Value = 179.5

$ws.Range("H58").Value = 649.5
$ws.Range("I58").Value = 400
$ws.Range("K58").Value = 1200
$ws.Range("M58").Value = -1072

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H7").Value = 1000000
$ws.Range("I7").Value = 1000000
$ws.Range("J7").Value = 0
$ws.Range("K7").Value = 1000000
$ws.Range("L7").ClearContents()
$ws.Range("M7").Value = -999888
$ws.Range("N7").Value = 0

$ws.Range("H8").Value = 1000000
$ws.Range("I8").Value = 1000000
$ws.Range("J8").Value = 0
$ws.Range("K8").Value = 1000000
$ws.Range("L8").ClearContents()
$ws.Range("M8").Value = -999861
$ws.Range("N8").Value = 0

$ws.Range("H10").Value = 40000000
$ws.Range("I10").Value = 40000000
$ws.Range("K10").Value = 40000000
$ws.Range("M10").Value = -39999831

$ws.Range("H11").Value = 32857428
$ws.Range("I11").Value = 38333332
$ws.Range("J11").Value = 2000
$ws.Range("K11").Value = 38333332
$ws.Range("L11").Value = 2000
$ws.Range("M11").Value = -38333193
$ws.Range("N11").Value = -2278

$ws.Range("H13").Value = 3600
$ws.Range("I13").Value = 3600
$ws.Range("K13").Value = 3600
$ws.Range("M13").Value = -3461

$ws.Range("H55").Value = 5250
$ws.Range("I55").Value = 3500
$ws.Range("K55").Value = 3500
$ws.Range("M55").Value = -3173

$ws.Range("H102").Value = 2590.8572
$ws.Range("I102").Value = 2269.4
$ws.Range("J102").Value = 3394.5
$ws.Range("K102").Value = 2269.4
$ws.Range("L102").Value = 3394.5
$ws.Range("M102").Value = -647.4000000000001
$ws.Range("N102").Value = -6638.5

$ws.Range("H113").Value = 19247352
$ws.Range("I113").Value = 27787678
$ws.Range("J113").Value = 31622.5
$ws.Range("K113").Value = 27787678
$ws.Range("L113").Value = 31622.5
$ws.Range("M113").Value = -27785508
$ws.Range("N113").Value = -35962.5

$ws.Range("H114").Value = 30000
$ws.Range("J114").Value = 30000
$ws.Range("L114").Value = 30000
$ws.Range("N114").Value = -38678

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H4").Value = 4983
$ws.Range("I4").Value = 6504.5
$ws.Range("J4").Value = 3968.6667
$ws.Range("K4").Value = 6504.5
$ws.Range("L4").Value = 3968.6667
$ws.Range("M4").Value = -6391.5
$ws.Range("N4").Value = -4194.6667

$ws.Range("H5").Value = 15011
$ws.Range("I5").Value = 0
$ws.Range("K5").Value = 0
$ws.Range("M5").ClearContents()

$ws.Range("H28").Value = 4983
$ws.Range("I28").Value = 6504.5
$ws.Range("J28").Value = 3968.6667
$ws.Range("K28").Value = 6504.5
$ws.Range("L28").Value = 3968.6667
$ws.Range("M28").Value = -6272.5
$ws.Range("N28").Value = -4432.6667

$ws.Range("H37").Value = 4983
$ws.Range("I37").Value = 6504.5
$ws.Range("J37").Value = 3968.6667
$ws.Range("K37").Value = 6504.5
$ws.Range("L37").Value = 3968.6667
$ws.Range("M37").Value = -6397.5
$ws.Range("N37").Value = -4182.6667

$ws.Range("H41").Value = 0
$ws.Range("I41").Value = 0
$ws.Range("J41").Value = 0
$ws.Range("K41").Value = 0
$ws.Range("L41").ClearContents()
$ws.Range("M41").ClearContents()
$ws.Range("N41").Value = 0

$ws.Range("H47").Value = 0
$ws.Range("I47").Value = 0
$ws.Range("K47").Value = 0
$ws.Range("M47").ClearContents()

$ws.Range("H52").Value = 0
$ws.Range("I52").Value = 0
$ws.Range("K52").Value = 0
$ws.Range("M52").ClearContents()

$ws.Range("H100").Value = 12498.5
$ws.Range("I100").Value = 12498.5
$ws.Range("J100").Value = 0
$ws.Range("K100").Value = 12498.5
$ws.Range("L100").Value = 0
$ws.Range("M100").ClearContents()
$ws.Range("N100").Value = -11957.5

$ws.Range("H122").Value = 1829.9
$ws.Range("I122").Value = 1724.875
$ws.Range("K122").Value = 5174.625
$ws.Range("M122").Value = -2724.625

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H12").Value = 2965.6667
$ws.Range("I12").Value = 2900
$ws.Range("J12").Value = 2998.5
$ws.Range("K12").Value = 2900
$ws.Range("L12").Value = 2998.5
$ws.Range("M12").Value = -2758
$ws.Range("N12").Value = -3282.5

$ws.Range("H29").Value = 10015000
$ws.Range("I29").Value = 10015000
$ws.Range("K29").Value = 10015000
$ws.Range("M29").Value = -10014710

$ws.Range("H99").Value = 25000
$ws.Range("I99").Value = 25000
$ws.Range("K99").Value = 25000
$ws.Range("M99").Value = -22005

$ws.Range("H117").Value = 75000
$ws.Range("J117").Value = 75000
$ws.Range("L117").Value = 75000
$ws.Range("N117").Value = -84178

$ws.Range("H130").Value = 75000
$ws.Range("J130").Value = 75000
$ws.Range("L130").Value = 75000
$ws.Range("N130").Value = -85040
